$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "2023" column (Q) is being appended after the existing "2022" column (P)
# for every data row of the table (header row 2, year-label row 3, and the
# three data rows 4-6).

# --- Row 2 (blank spacer/border row under the title) ---
$ws.Range("Q2").Value = $null
$ws.Range("P2").Copy()
$ws.Range("Q2").PasteSpecial(-4122)   # xlPasteFormats - copy P2's formatting onto Q2

# --- Row 3 (year headers) ---
$ws.Range("Q3").Value = 2023
$ws.Range("P3").Copy()
$ws.Range("Q3").PasteSpecial(-4122)

# --- Row 4 (per-person waste, kg) ---
$ws.Range("Q4").Value = 279.01945525291825
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)

# --- Row 5 (waste removed, thousand tons) ---
$ws.Range("Q5").Value = 1792.7
$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)

# --- Row 6 (average annual population, thousand people) ---
$ws.Range("Q6").Value = 6425
$ws.Range("P6").Copy()
$ws.Range("Q6").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Rows 4 and 5 grew slightly taller to fit the wrapped header text now that
# an extra column has been added.
$ws.Rows.Item(4).RowHeight = 27
$ws.Rows.Item(5).RowHeight = 27.75

# Reset the sheet's selection back to the top-left cell (the saved workbook
# previously had a stray selection sitting on S4).
$ws.Range("A1").Select() | Out-Null
